# Added term accession numbers to ambiguous tags and harmonized similar tags.
# Target sheet: SwateTemplateMetadata
#   Row 8  -> ER                         = PRIDE   (already set)
#   Row 9  -> ER Term Accession Number   -> http://purl.obolibrary.org/obo/DPBO_1000098
#   Row 10 -> ER Term Source REF         -> DPBO
#   Row 12 -> Tags (headers): Assay | Proteomics | Measurement | Mass spectrometry | MS | PRIDE
#   Row 13 -> Tags Term Accession Number -> C13 (Proteomics) = http://purl.obolibrary.org/obo/NCIT_C20085
#                                            G13 (PRIDE)      = http://purl.obolibrary.org/obo/DPBO_1000098
#   Row 14 -> Tags Term Source REF       -> C14 (Proteomics) = NCIT
#                                            G14 (PRIDE)      = DPBO

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SwateTemplateMetadata")

$ws.Range("C13").Value = "http://purl.obolibrary.org/obo/NCIT_C20085"
$ws.Range("C14").Value = "NCIT"

$ws.Range("B9").Value = "http://purl.obolibrary.org/obo/DPBO_1000098"
$ws.Range("G13").Value = "http://purl.obolibrary.org/obo/DPBO_1000098"

$ws.Range("B10").Value = "DPBO"
$ws.Range("G14").Value = "DPBO"

# Row 13 now wraps the long accession-number URL in C13, so its height grows
# to the value Excel computed on re-save.
$ws.Rows.Item(13).RowHeight = 57.6

# Move the active selection/cursor on this sheet from G13 to B13, matching
# the authored workbook state after the edit.
$ws.Activate()
$ws.Range("B13").Select()
